# Add "N. crash:", "Mean:" and "Standard Deviation:" summary rows below the
# data table (rows 53-55), matching the commit "Add Mean and Standard Deviation".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 53: crash counter
$ws.Range("F53").Value = "N. crash:"
$ws.Range("G53").Formula = '=COUNTIF(G2:G51,"True")'

# Row 54: column-wise mean (I:P)
$ws.Range("F54").Value = "Mean:"
$ws.Range("I54:P54").Formula = "=AVERAGE(I2:I51)"

# Row 55: column-wise sample standard deviation (I:P)
$ws.Range("F55").Value = "Standard Deviation:"
$ws.Range("I55:P55").Formula = "=STDEV.S(I2:I51)"

# Move the active selection near the new rows, like the authored workbook.
[void]$ws.Range("K60").Select()
